$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C39").Value = 'Arica y Parinacota'
$ws.Range("D39").Value = 44641
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 100112031
$ws.Range("G39").Value = 'Poroto verde'
$ws.Range("H39").Value = 'Sin especificar'
$ws.Range("I39").Value = 'Primera'
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 1300
$ws.Range("L39").Value = 1400
$ws.Range("M39").Value = 1350
$ws.Range("N39").Value = '$/kilo'
$ws.Range("O39").Value = 'Región de Arica y Parinacota'
$ws.Range("P39").Value = 1350
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = 'Hortaliza'

# Row 40
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C40").Value = 'Arica y Parinacota'
$ws.Range("D40").Value = 44483
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112031
$ws.Range("G40").Value = 'Poroto verde'
$ws.Range("H40").Value = 'Sin especificar'
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 1100
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = 1250
$ws.Range("N40").Value = '$/kilo'
$ws.Range("O40").Value = 'Región de Arica y Parinacota'
$ws.Range("P40").Value = 1250
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = 'Hortaliza'

# Row 41
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C41").Value = 'Arica y Parinacota'
$ws.Range("D41").Value = 44264
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112031
$ws.Range("G41").Value = 'Poroto verde'
$ws.Range("H41").Value = 'Sin especificar'
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 1300
$ws.Range("K41").Value = 2200
$ws.Range("L41").Value = 2300
$ws.Range("M41").Value = 2250
$ws.Range("N41").Value = '$/kilo'
$ws.Range("O41").Value = 'Región de Arica y Parinacota'
$ws.Range("P41").Value = 2250
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = 'Hortaliza'

# Row 42
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C42").Value = 'Arica y Parinacota'
$ws.Range("D42").Value = 44330
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = 100112031
$ws.Range("G42").Value = 'Poroto verde'
$ws.Range("H42").Value = 'Magnum'
$ws.Range("I42").Value = 'Primera'
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 24000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 24500
$ws.Range("N42").Value = '$/malla 25 kilos'
$ws.Range("O42").Value = 'Perú'
$ws.Range("P42").Value = 980
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = 'Hortaliza'

# Row 43
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C43").Value = 'Arica y Parinacota'
$ws.Range("D43").Value = 44370
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = 100112031
$ws.Range("G43").Value = 'Poroto verde'
$ws.Range("H43").Value = 'Magnum'
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 80
$ws.Range("K43").Value = 19000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 19375
$ws.Range("N43").Value = '$/malla 25 kilos'
$ws.Range("O43").Value = 'Perú'
$ws.Range("P43").Value = 775
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = 'Hortaliza'

# Row 44
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C44").Value = 'Arica y Parinacota'
$ws.Range("D44").Value = 44370
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100112031
$ws.Range("G44").Value = 'Poroto verde'
$ws.Range("H44").Value = 'Magnum'
$ws.Range("I44").Value = 'Segunda'
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 17000
$ws.Range("L44").Value = 18000
$ws.Range("M44").Value = 17500
$ws.Range("N44").Value = '$/malla 25 kilos'
$ws.Range("O44").Value = 'Perú'
$ws.Range("P44").Value = 700
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = 'Hortaliza'

# Row 45
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C45").Value = 'Arica y Parinacota'
$ws.Range("D45").Value = 44385
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112031
$ws.Range("G45").Value = 'Poroto verde'
$ws.Range("H45").Value = 'Sin especificar'
$ws.Range("I45").Value = 'Primera'
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1200
$ws.Range("L45").Value = 1300
$ws.Range("M45").Value = 1250
$ws.Range("N45").Value = '$/kilo'
$ws.Range("O45").Value = 'Región de Arica y Parinacota'
$ws.Range("P45").Value = 1250
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = 'Hortaliza'

# Row 46
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C46").Value = 'Arica y Parinacota'
$ws.Range("D46").Value = 44312
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100112031
$ws.Range("G46").Value = 'Poroto verde'
$ws.Range("H46").Value = 'Sin especificar'
$ws.Range("I46").Value = 'Primera'
$ws.Range("J46").Value = 1700
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 1400
$ws.Range("M46").Value = 1350
$ws.Range("N46").Value = '$/kilo'
$ws.Range("O46").Value = 'Región de Arica y Parinacota'
$ws.Range("P46").Value = 1350
$ws.Range("Q46").Value = 1
$ws.Range("R46").Value = 'Hortaliza'

# Row 47
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C47").Value = 'Arica y Parinacota'
$ws.Range("D47").Value = 44399
$ws.Range("E47").Value = 15
$ws.Range("F47").Value = 100112031
$ws.Range("G47").Value = 'Poroto verde'
$ws.Range("H47").Value = 'Magnum'
$ws.Range("I47").Value = 'Primera'
$ws.Range("J47").Value = 1400
$ws.Range("K47").Value = 1300
$ws.Range("L47").Value = 1400
$ws.Range("M47").Value = 1350
$ws.Range("N47").Value = '$/kilo'
$ws.Range("O47").Value = 'Región de Arica y Parinacota'
$ws.Range("P47").Value = 1350
$ws.Range("Q47").Value = 1
$ws.Range("R47").Value = 'Hortaliza'

# Row 48
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C48").Value = 'Arica y Parinacota'
$ws.Range("D48").Value = 44543
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100112031
$ws.Range("G48").Value = 'Poroto verde'
$ws.Range("H48").Value = 'Sin especificar'
$ws.Range("I48").Value = 'Primera'
$ws.Range("J48").Value = 1300
$ws.Range("K48").Value = 400
$ws.Range("L48").Value = 450
$ws.Range("M48").Value = 425
$ws.Range("N48").Value = '$/kilo'
$ws.Range("O48").Value = 'Región de Arica y Parinacota'
$ws.Range("P48").Value = 425
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = 'Hortaliza'

# Row 49
$ws.Range("A49").Value = 1
$ws.Range("B49").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C49").Value = 'Arica y Parinacota'
$ws.Range("D49").Value = 44258
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = 100112031
$ws.Range("G49").Value = 'Poroto verde'
$ws.Range("H49").Value = 'Sin especificar'
$ws.Range("I49").Value = 'Primera'
$ws.Range("J49").Value = 1600
$ws.Range("K49").Value = 2300
$ws.Range("L49").Value = 2500
$ws.Range("M49").Value = 2400
$ws.Range("N49").Value = '$/kilo'
$ws.Range("O49").Value = 'Región de Arica y Parinacota'
$ws.Range("P49").Value = 2400
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = 'Hortaliza'

# Row 50
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C50").Value = 'Arica y Parinacota'
$ws.Range("D50").Value = 44390
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = 'Poroto verde'
$ws.Range("H50").Value = 'Magnum'
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 1500
$ws.Range("K50").Value = 900
$ws.Range("L50").Value = 1000
$ws.Range("M50").Value = 950
$ws.Range("N50").Value = '$/kilo'
$ws.Range("O50").Value = 'Región de Arica y Parinacota'
$ws.Range("P50").Value = 950
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = 'Hortaliza'

# Row 51
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C51").Value = 'Arica y Parinacota'
$ws.Range("D51").Value = 44179
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = 100112031
$ws.Range("G51").Value = 'Poroto verde'
$ws.Range("H51").Value = 'Sin especificar'
$ws.Range("I51").Value = 'Primera'
$ws.Range("J51").Value = 1600
$ws.Range("K51").Value = 1200
$ws.Range("L51").Value = 1300
$ws.Range("M51").Value = 1250
$ws.Range("N51").Value = '$/kilo'
$ws.Range("O51").Value = 'Región de Arica y Parinacota'
$ws.Range("P51").Value = 1250
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = 'Hortaliza'

# Row 52
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C52").Value = 'Arica y Parinacota'
$ws.Range("D52").Value = 44160
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = 100112031
$ws.Range("G52").Value = 'Poroto verde'
$ws.Range("H52").Value = 'Sin especificar'
$ws.Range("I52").Value = 'Primera'
$ws.Range("J52").Value = 1900
$ws.Range("K52").Value = 700
$ws.Range("L52").Value = 800
$ws.Range("M52").Value = 750
$ws.Range("N52").Value = '$/kilo'
$ws.Range("O52").Value = 'Región de Arica y Parinacota'
$ws.Range("P52").Value = 750
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = 'Hortaliza'
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
